$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 100 - this pushes the existing rows
# 100-105 down to 102-107, preserving all their data/formatting.
$ws.Rows("100:101").Insert()

# Row 100: new weekly entry for "Argentina(o)" variety
$ws.Range("A100").Value = 11
$ws.Range("B100").Value = "Vega Monumental Concepción"
$ws.Range("C100").Value = "Bíobío"
$ws.Range("D100").Value = 45147
$ws.Range("D100").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E100").Value = 8
$ws.Range("F100").Value = 100112013
$ws.Range("G100").Value = "Alcachofa"
$ws.Range("H100").Value = "Argentina(o)"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 50
$ws.Range("K100").Value = 13000
$ws.Range("L100").Value = 13000
$ws.Range("M100").Value = 13000
$ws.Range("N100").Value = "`$/caja 50 unidades"
$ws.Range("O100").Value = "Provincia de Limarí"
$ws.Range("P100").Value = 260
$ws.Range("Q100").Value = 50
$ws.Range("R100").Value = "Hortaliza"

# Row 101: new weekly entry for "Española" variety
$ws.Range("A101").Value = 11
$ws.Range("B101").Value = "Vega Monumental Concepción"
$ws.Range("C101").Value = "Bíobío"
$ws.Range("D101").Value = 45147
$ws.Range("D101").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E101").Value = 8
$ws.Range("F101").Value = 100112013
$ws.Range("G101").Value = "Alcachofa"
$ws.Range("H101").Value = "Española"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 50
$ws.Range("K101").Value = 13000
$ws.Range("L101").Value = 13000
$ws.Range("M101").Value = 13000
$ws.Range("N101").Value = "`$/caja 30 unidades"
$ws.Range("O101").Value = "Provincia de Limarí"
$ws.Range("P101").Value = 433
$ws.Range("Q101").Value = 30
$ws.Range("R101").Value = "Hortaliza"
